$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Row 1 header: F1 ("Moisture Content") gets a new style
#    (wrap text + left-aligned) instead of the plain header style.
# ---------------------------------------------------------------------
$ws.Range("F1").HorizontalAlignment = -4131   # xlLeft
$ws.Range("F1").WrapText = $true

# ---------------------------------------------------------------------
# 2. Columns E (5) / F (6) get slightly wider, and column E loses its
#    "best fit" flag (we just set an explicit width for both).
# ---------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 7.666666666666667
$ws.Columns.Item(6).ColumnWidth = 7.830729166666667

# ---------------------------------------------------------------------
# 3. The old row 10 ("electricity-eurofer") is removed entirely, which
#    shifts the old rows 11-13 (coke / coking coal / PCI coal) up to
#    10-12.
# ---------------------------------------------------------------------
$ws.Rows.Item(10).Delete()

# ---------------------------------------------------------------------
# 4. Append the new Eurofer low-carbon-roadmap "PROXY fuel mix" rows
#    13-15 (2010 / 2030 / 2050).
# ---------------------------------------------------------------------

# Row 13 - PROXY fuel mix (1:1 energy:mass unit)- Eurofer Electricity 2010
$ws.Range("A13").Value = "PROXY fuel mix (1:1 energy:mass unit)- Eurofer Electricity 2010"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 0.11
$ws.Range("E13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 1
$ws.Range("O13").Formula = "=1-N13"
$ws.Range("Q13").Value = "EUROFER"

# Row 14 - PROXY fuel mix (1:1 energy:mass unit)- Eurofer Electricity 2030
$ws.Range("A14").Value = "PROXY fuel mix (1:1 energy:mass unit)- Eurofer Electricity 2030"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 0.055
$ws.Range("E14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 1
$ws.Range("O14").Formula = "=1-N14"
$ws.Range("Q14").Value = "EUROFER"

# Row 15 - PROXY fuel mix (1:1 energy:mass unit)- Eurofer Electricity 2050
$ws.Range("A15").Value = "PROXY fuel mix (1:1 energy:mass unit)- Eurofer Electricity 2050"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 1
$ws.Range("O15").Formula = "=1-N15"
$ws.Range("Q15").Value = "EUROFER"

# ---------------------------------------------------------------------
# 5. Restore the selected cell to match the saved view state.
# ---------------------------------------------------------------------
$ws.Range("Q8").Select() | Out-Null
